$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122, xlPasteValues = -4163
$xlPasteFormats = -4122
$xlPasteValues = -4163

# ---------------------------------------------------------------
# Shift the trailing rows down by one to make room for the new
# "STERONATE  5MG 20 TAB" line item that was sold, inserted as the
# new row 13. Work bottom-up so we never clobber data we still need.
# ---------------------------------------------------------------

# Footer row (row number / "developed by") 17 -> 18
$ws.Range("A17:Q17").Copy()
$ws.Range("A18:Q18").PasteSpecial($xlPasteFormats)
$ws.Range("A17:Q17").Copy()
$ws.Range("A18:Q18").PasteSpecial($xlPasteValues)
$ws.Range("A18:F18").Merge()
$ws.Range("G18:I18").Merge()
$ws.Range("K18:Q18").Merge()

# Totals row 16 -> 17
$ws.Range("A16:Q16").Copy()
$ws.Range("A17:Q17").PasteSpecial($xlPasteFormats)
$ws.Range("A16:Q16").Copy()
$ws.Range("A17:Q17").PasteSpecial($xlPasteValues)
$ws.Range("P17:Q17").Merge()
$ws.Range("P17").Value = 515.5

# Old row 15 ("معجون كلوز اب الصغير") -> row 16, now a regular item row
$ws.Range("A15:Q15").Copy()
$ws.Range("A16:Q16").PasteSpecial($xlPasteFormats)
$ws.Range("A15:Q15").Copy()
$ws.Range("A16:Q16").PasteSpecial($xlPasteValues)
$ws.Range("A16:B16").Merge()
$ws.Range("C16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()
$ws.Range("N16:O16").Merge()
$excel.CutCopyMode = 0

$ws.Range("A16").Value = 10
$ws.Range("C16").Value = "معجون كلوز اب الصغير"
$ws.Range("H16").Value = "4:0"
$ws.Range("L16").Value = "0"
$ws.Range("N16").Value = "20.00"
$ws.Range("P16").Value = "20.0000"
$ws.Range("Q16").Value = "1:0"

# Old row 14 ("سرنجات 5 سم") -> row 15 (formatting/merges already correct there)
$ws.Range("A15").Value = 9
$ws.Range("C15").Value = "سرنجات 5 سم"
$ws.Range("H15").Value = "0:0"
$ws.Range("L15").Value = "0"
$ws.Range("N15").Value = "3.00"
$ws.Range("P15").Value = "3.0000"
$ws.Range("Q15").Value = "1:0"

# Old row 13 ("ديتول صغير") -> row 14 (formatting/merges already correct there)
$ws.Range("A14").Value = 8
$ws.Range("C14").Value = "ديتول صغير"
$ws.Range("H14").Value = "3:0"
$ws.Range("L14").Value = "0"
$ws.Range("N14").Value = "17.00"
$ws.Range("P14").Value = "17.0000"
$ws.Range("Q14").Value = "1:0"

# New row 13: "STERONATE  5MG 20 TAB" (formatting/merges already correct there)
$ws.Range("A13").Value = 7
$ws.Range("C13").Value = "STERONATE  5MG 20 TAB"
$ws.Range("H13").Value = "-1:0"
$ws.Range("L13").Value = "1"
$ws.Range("N13").Value = "56.00"
$ws.Range("P13").Value = "112.0000"
$ws.Range("Q13").Value = "2:0"

# Update the printed timestamp in the footer to reflect the re-upload time
$ws.Range("A18").Value = "Monday, 11 August, 2025 11:14 AM"
